# Adds new testscript WAT40 to the "Test Cases" sheet (sheet1), as row 27.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 26 is the last existing data row; copy its A:D formatting down to the
# new row 27 so the new row matches the look of the existing table rows.
$ws.Range("A26:D26").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new test case values.
$ws.Range("A27").Value = "WAT40"
$ws.Range("B27").Value = "WAT-147"
$ws.Range("C27").Value = "Verify that FIND button will be disabled at the beginning."
$ws.Range("D27").Value = "Y"

# Leave the user's selection on the newly added row, matching the authored
# edit (selection moved to C27 after entering the new test case).
$ws.Range("C27").Select() | Out-Null
